# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the a3c446ec-2b83-4084-9dfd-c83f819c9d91
# report row across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for a3c446ec-... row (row 2)
$overview.Range("G2").Value = "2016-09-02 19:14:13"

# zh-cn sheet, row 2 (a3c446ec-... handoff/handback)
$zhcn.Range("H2").Value = "2016-09-02 19:14:03"
$zhcn.Range("K2").Value = "2016-09-02 19:14:31"

# de-de sheet, row 2 (a3c446ec-... handoff/handback)
# Note: H2 shares the same text as Overview!G2 ("Latest HO Xliff Generate Date"
# value), so it must be updated to the same new timestamp as well.
$dede.Range("H2").Value = "2016-09-02 19:14:13"
$dede.Range("K2").Value = "2016-09-02 19:14:39"
